# DOVPROEV-6902 finalisatie eerste publicatie
#
# The 'bekistingsmaterialen' concept scheme is finalized for its first
# publication: two new columns ('definition' and 'note') are inserted
# after 'prefLabel' (pushing topConceptOf/hasTopConcept from H:I to J:K),
# and the concept rows are replaced by the finalized set of 5 concepts
# (notations 0-4) with their labels, definitions and notes, plus the
# collection and concept-scheme summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at H:I (definition, note) -- shifts the old
# topConceptOf/hasTopConcept columns from H:I to J:K.
$ws.Range("H1:I1").EntireColumn.Insert()

# Drop the now-obsolete rows 9-16 (the old, non-finalized concept list);
# the finalized sheet only spans 8 rows (header + 7 data rows).
$ws.Rows.Item(9).Resize(8).Delete()

# Keep the numeric-looking notations ('0'..'4') stored as text, matching
# the source data (avoids Excel auto-converting them to numbers).
$ws.Range("F3:F7").NumberFormat = "@"

# Write out the finalized A1:K8 content cell by cell (multi-cell array
# assignment isn't reliable against this host, so address each cell
# individually via Cells.Item).
$ws.Cells.Item(1,1).Value = 'id'  # A1
$ws.Cells.Item(1,2).Value = 'type'  # B1
$ws.Cells.Item(1,3).Value = 'dc\.identifier'  # C1
$ws.Cells.Item(1,4).Value = 'inScheme'  # D1
$ws.Cells.Item(1,5).Value = 'member'  # E1
$ws.Cells.Item(1,6).Value = 'notation'  # F1
$ws.Cells.Item(1,7).Value = 'prefLabel'  # G1
$ws.Cells.Item(1,8).Value = 'definition'  # H1
$ws.Cells.Item(1,9).Value = 'note'  # I1
$ws.Cells.Item(1,10).Value = 'topConceptOf'  # J1
$ws.Cells.Item(1,11).Value = 'hasTopConcept'  # K1

$ws.Cells.Item(2,1).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/collection/bekistingsmateriaal/bekistingsmaterialen'  # A2
$ws.Cells.Item(2,2).Value = 'http://www.w3.org/2004/02/skos/core#Collection'  # B2
$ws.Cells.Item(2,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.collection.bekistingsmateriaal.bekistingsmaterialen'  # C2
$ws.Cells.Item(2,4).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # D2
$ws.Cells.Item(2,5).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/0|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/1|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/2|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/3|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/4'  # E2
$ws.Cells.Item(2,6).Value = 'collectie_bekistingsmaterialen'  # F2
$ws.Cells.Item(2,7).Value = 'Collectie van bekistingsmaterialen.'  # G2
$ws.Cells.Item(2,8).Value = 'null'  # H2
$ws.Cells.Item(2,9).Value = 'null'  # I2
$ws.Cells.Item(2,10).Value = 'null'  # J2
$ws.Cells.Item(2,11).Value = 'null'  # K2

$ws.Cells.Item(3,1).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/0'  # A3
$ws.Cells.Item(3,2).Value = 'http://www.w3.org/2004/02/skos/core#Concept'  # B3
$ws.Cells.Item(3,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.bekistingsmateriaal.0'  # C3
$ws.Cells.Item(3,4).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # D3
$ws.Cells.Item(3,5).Value = 'null'  # E3
$ws.Cells.Item(3,6).Value = '0'  # F3
$ws.Cells.Item(3,7).Value = 'natuurlijke materialen (biodegradeerbaar, hout, …)'  # G3
$ws.Cells.Item(3,8).Value = 'Het boorgat is afgewerkt met natuurlijke materialen zoals bijvoorbeeld biologisch afbreekbare stoffen, hout.'  # H3
$ws.Cells.Item(3,9).Value = 'Het boorgat is afgewerkt met natuurlijke materialen zoals bijvoorbeeld biologisch afbreekbare stoffen, hout.'  # I3
$ws.Cells.Item(3,10).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # J3
$ws.Cells.Item(3,11).Value = 'null'  # K3

$ws.Cells.Item(4,1).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/1'  # A4
$ws.Cells.Item(4,2).Value = 'http://www.w3.org/2004/02/skos/core#Concept'  # B4
$ws.Cells.Item(4,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.bekistingsmateriaal.1'  # C4
$ws.Cells.Item(4,4).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # D4
$ws.Cells.Item(4,5).Value = 'null'  # E4
$ws.Cells.Item(4,6).Value = '1'  # F4
$ws.Cells.Item(4,7).Value = 'kunststof: PVC'  # G4
$ws.Cells.Item(4,8).Value = 'De materialen die in het boorgat achterblijven als al dan niet verloren bekisting of buizen bestaan uit PVC.'  # H4
$ws.Cells.Item(4,9).Value = 'De materialen die in het boorgat achterblijven als al dan niet verloren bekisting of buizen bestaan uit PVC.'  # I4
$ws.Cells.Item(4,10).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # J4
$ws.Cells.Item(4,11).Value = 'null'  # K4

$ws.Cells.Item(5,1).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/2'  # A5
$ws.Cells.Item(5,2).Value = 'http://www.w3.org/2004/02/skos/core#Concept'  # B5
$ws.Cells.Item(5,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.bekistingsmateriaal.2'  # C5
$ws.Cells.Item(5,4).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # D5
$ws.Cells.Item(5,5).Value = 'null'  # E5
$ws.Cells.Item(5,6).Value = '2'  # F5
$ws.Cells.Item(5,7).Value = 'kunststof: andere (PE, HDPE, LDPE, …)'  # G5
$ws.Cells.Item(5,8).Value = 'De materialen die in het boorgat achterblijven als al dan niet verloren bekisting of buizen bestaan uit kunststof anders dan PVC zoals polyethylene (PE, HDPE, LDPE, …).'  # H5
$ws.Cells.Item(5,9).Value = 'De materialen die in het boorgat achterblijven als al dan niet verloren bekisting of buizen bestaan uit kunststof anders dan PVC zoals polyethylene (PE, HDPE, LDPE, …).'  # I5
$ws.Cells.Item(5,10).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # J5
$ws.Cells.Item(5,11).Value = 'null'  # K5

$ws.Cells.Item(6,1).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/3'  # A6
$ws.Cells.Item(6,2).Value = 'http://www.w3.org/2004/02/skos/core#Concept'  # B6
$ws.Cells.Item(6,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.bekistingsmateriaal.3'  # C6
$ws.Cells.Item(6,4).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # D6
$ws.Cells.Item(6,5).Value = 'null'  # E6
$ws.Cells.Item(6,6).Value = '3'  # F6
$ws.Cells.Item(6,7).Value = 'metaal: inox, staal, roestvrij staal, gegalvaniseerd staal, …'  # G6
$ws.Cells.Item(6,8).Value = 'De materialen die in het boorgat achterblijven als al dan niet verloren bekisting of buizen bestaan uit metaal zoals inox of staal.'  # H6
$ws.Cells.Item(6,9).Value = 'De materialen die in het boorgat achterblijven als al dan niet verloren bekisting of buizen bestaan uit metaal zoals inox of staal.'  # I6
$ws.Cells.Item(6,10).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # J6
$ws.Cells.Item(6,11).Value = 'null'  # K6

$ws.Cells.Item(7,1).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/4'  # A7
$ws.Cells.Item(7,2).Value = 'http://www.w3.org/2004/02/skos/core#Concept'  # B7
$ws.Cells.Item(7,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.bekistingsmateriaal.4'  # C7
$ws.Cells.Item(7,4).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # D7
$ws.Cells.Item(7,5).Value = 'null'  # E7
$ws.Cells.Item(7,6).Value = '4'  # F7
$ws.Cells.Item(7,7).Value = 'beton/metselwerk'  # G7
$ws.Cells.Item(7,8).Value = 'Het boorgat is afgewerkt met bekisting uit beton of metselwerk.'  # H7
$ws.Cells.Item(7,9).Value = 'Het boorgat is afgewerkt met bekisting uit beton of metselwerk.'  # I7
$ws.Cells.Item(7,10).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # J7
$ws.Cells.Item(7,11).Value = 'null'  # K7

$ws.Cells.Item(8,1).Value = 'https://data.omgeving.vlaanderen.be/id/conceptscheme/bekistingsmateriaal'  # A8
$ws.Cells.Item(8,2).Value = 'http://www.w3.org/2004/02/skos/core#ConceptScheme'  # B8
$ws.Cells.Item(8,3).Value = 'be.vlaanderen.bodemenondergrond.data.id.conceptscheme.bekistingsmateriaal'  # C8
$ws.Cells.Item(8,4).Value = 'null'  # D8
$ws.Cells.Item(8,5).Value = 'null'  # E8
$ws.Cells.Item(8,6).Value = 'bekistingsmaterialen'  # F8
$ws.Cells.Item(8,7).Value = 'Conceptschema van bekistingsmaterialen.'  # G8
$ws.Cells.Item(8,8).Value = 'null'  # H8
$ws.Cells.Item(8,9).Value = 'null'  # I8
$ws.Cells.Item(8,10).Value = 'null'  # J8
$ws.Cells.Item(8,11).Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/0|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/1|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/2|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/3|https://data.bodemenondergrond.vlaanderen.be/id/concept/bekistingsmateriaal/4'  # K8

Write-Host "done"
